$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are free-form text (e.g. "65.482.40" thousand-grouped
# numbers, or plain decimals). Force text storage via NumberFormat "@" so that
# Excel does not reinterpret values such as "595.02" as numbers.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.482.40'
$ws.Range('E2').Value = '  -1.95%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.403.62'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '595.02'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.37'
$ws.Range('E6').Value = '  -4.03%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.403.32'
$ws.Range('E7').Value = '  -1.66%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.469'
$ws.Range('E9').Value = '  -2.98%  '
$ws.Range('E10').Value = '  -5.18%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.87'
$ws.Range('E11').Value = '  +5.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.404'
$ws.Range('E12').Value = '  -4.59%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.975.62'
$ws.Range('E13').Value = '  -1.88%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000200'
$ws.Range('E14').Value = '  -6.46%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '29.71'
$ws.Range('E15').Value = '  -6.07%  '
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.395.61'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '65.454.22'
$ws.Range('E18').Value = '  -2.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.46'
$ws.Range('E19').Value = '  +4.58%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.11'
$ws.Range('E20').Value = '  -5.10%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.75'
$ws.Range('E21').Value = '  -3.71%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '416.53'
$ws.Range('E22').Value = '  -5.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.580'
$ws.Range('E23').Value = '  -4.77%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '77.38'
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.539.14'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000111'
$ws.Range('E27').Value = '  -8.41%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.24'
$ws.Range('E28').Value = '  -6.62%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.84'
$ws.Range('E29').Value = '  -6.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.42'
$ws.Range('E30').Value = '  -2.57%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.161'
$ws.Range('E32').Value = '  -3.73%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.47'
$ws.Range('E33').Value = '  -8.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '24.49'
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.394.84'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.70'
$ws.Range('E37').Value = '  -5.81%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.54'
$ws.Range('E38').Value = '  -8.80%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.55'
$ws.Range('E39').Value = '  -4.56%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '170.06'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0859'
$ws.Range('E42').Value = '  -3.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.06'
$ws.Range('E43').Value = '  -6.01%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.870'
$ws.Range('E44').Value = '  -1.44%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.92'
$ws.Range('E45').Value = '  -11.54%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '45.50'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.85'
$ws.Range('E47').Value = '  -7.36%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.17'
$ws.Range('E48').Value = '  -5.70%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.10'
$ws.Range('E49').Value = '  -4.92%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.28'
$ws.Range('E50').Value = '  -6.96%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.919'
$ws.Range('E51').Value = '  -6.51%  '
